# Update Work Week and Social Spending
# (Armenia GDP per Capita data: revise existing series 1973/1980-2010 and
#  append newly available years 2011-2016.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Helper: write a value into column E as TEXT (matches existing column E
# cells, which are all shared-string / text typed, e.g. "6152", "9806", ...)
# without leaving a residual number-format style on the cell.
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# --- Revise existing "Data" values in column E (rows 2, 9-39) ---------------

$updates = @{
    2  = "9806"
    9  = "9291"
    10 = "9548"
    11 = "9809"
    12 = "9865"
    13 = "10181"
    14 = "10485"
    15 = "10394"
    16 = "9789"
    17 = "9301"
    18 = "10525"
    19 = "9669"
    20 = "8141.67112356872"
    21 = "4699.61095270265"
    22 = "4130.1472172755"
    23 = "4392.97987658491"
    24 = "4703.13073351788"
    25 = "4811.48954701875"
    26 = "4836.04422056585"
    27 = "4994.05303977392"
    28 = "5000.45244935948"
    29 = "5139.82550006292"
    30 = "5457.5100620274"
    31 = "6082.14618880966"
    32 = "6743.24097825047"
    33 = "7229.6482263198"
    34 = "8006.96042553147"
    35 = "8798.75018914713"
    36 = "9713.21328027296"
    37 = "10080.8128333688"
    38 = "8399.08687296095"
    39 = "8330.81117545937"
}

foreach ($row in $updates.Keys) {
    Set-TextValue $row 5 $updates[$row]
}

# --- Append new "Data" rows for years 2011-2016 (rows 40-45) ----------------

$newRows = @(
    @{ Row = 40; Year = 2011; Value = "8465" }
    @{ Row = 41; Year = 2012; Value = "9077" }
    @{ Row = 42; Year = 2013; Value = "9385" }
    @{ Row = 43; Year = 2014; Value = "9735" }
    @{ Row = 44; Year = 2015; Value = "10042" }
    @{ Row = 45; Year = 2016; Value = "10080" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 51
    $ws.Cells.Item($row, 2).Value = "Armenia"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $r.Year
    Set-TextValue $row 5 $r.Value
}
